# Loose Coupling - Reflection concept
# Collapse the two-sheet "Negative Scenario" / "Positive Scenario" login
# test-data workbook into a single consolidated "Login Scenario" sheet
# with a Test Case Name / username / password table.

$wb = $excel.ActiveWorkbook

# --- Remove the "Positive Scenario" sheet, keep & repurpose the first one ---
$wb.Worksheets("Positive Scenario").Delete() | Out-Null

$ws = $wb.Worksheets("Negative Scenario")
$ws.Name = "Login Scenario"

# --- Column widths (approximate best-fit target widths) ---
$ws.Columns("A").ColumnWidth = 41.6
$ws.Columns("B").ColumnWidth = 15.4
$ws.Columns("C").ColumnWidth = 17.4

# --- Header row ---
$ws.Range("A1").Value = "Test Case Name"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "password"

# --- Data rows ---
$ws.Range("A2").Value = "Enter invalid username and invalid password"
$ws.Range("B2").Value = "invalid_username"
$ws.Range("C2").Value = "invalid_password"

$ws.Range("A3").Value = "Enter invalid username and valid password"
$ws.Range("B3").Value = "invalid_username"
$ws.Range("C3").Value = "secret_sauce"

$ws.Range("A4").Value = "Enter valid username and invalid password"
$ws.Range("B4").Value = "standard_user"
$ws.Range("C4").Value = "invalid_password"

$ws.Range("A5").Value = "Enter null username and invalid password"
$ws.Range("B5").Value = "'"
$ws.Range("C5").Value = "invalid_password"

$ws.Range("A6").Value = "Enter valid username and null password"
$ws.Range("B6").Value = "standard_user"
$ws.Range("C6").Value = "'"

$ws.Range("A7").Value = "Enter null username and null password"
$ws.Range("B7").Value = "'"
$ws.Range("C7").Value = "'"

$ws.Range("A8").Value = "Enter valid username and valid password"
$ws.Range("B8").Value = "standard_user"
$ws.Range("C8").Value = "secret_sauce"

# --- Selection moves to D1, matching the saved view state ---
$ws.Range("D1").Select() | Out-Null
